# Apply the "automatic update" changes to the OKÄNT (unknown) log sheet:
#  - Column C ("Förändrad" / last-changed date) bumps from 45667 to 45668
#    for every data row (rows 2-36).
#  - Rows 35 and 36 swap their "Beteckning" (A) and "Area (ha)" (G) values,
#    i.e. row 35 becomes "A 60501-2024" / 0.6 ha and row 36 becomes
#    "A 60500-2024" / 0.8 ha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for all data rows (2 through 36).
for ($row = 2; $row -le 36; $row++) {
    $ws.Cells.Item($row, 3).Value = 45668
}

# Swap the Beteckning (A) and Area (G) values between row 35 and row 36.
$a35 = $ws.Cells.Item(35, 1).Value()
$a36 = $ws.Cells.Item(36, 1).Value()
$g35 = $ws.Cells.Item(35, 7).Value()
$g36 = $ws.Cells.Item(36, 7).Value()

$ws.Cells.Item(35, 1).Value = $a36
$ws.Cells.Item(36, 1).Value = $a35

$ws.Cells.Item(35, 7).Value = $g36
$ws.Cells.Item(36, 7).Value = $g35
